$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$samplesQuery = @'
MATCH (ss:study_subject)
WITH COLLECT(ss.study_subject_id) AS all_subjects
MATCH (samp:sample)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (samp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
 WHERE ss.disease_subtype IN ["Infiltrating Lobular Carcinoma"]
WITH
    distinct lp,
    toInteger(split(ss.study_subject_id,'-')[2]) AS subject_id_num,
    collect(distinct f.file_id) AS files,
    samp, ss, s, p, all_subjects
RETURN
 samp.sample_id AS `Sample ID`,
            ss.study_subject_id AS `Case ID`,
            p.program_acronym AS `Program Code`,
            s.study_acronym AS `Arm`,
            ss.disease_subtype AS `Diagnosis`,
            samp.tissue_type AS `Tissue Type`,
            samp.composition AS `Tissue Composition`,
            samp.sample_anatomic_site AS `Sample Anatomic Site`,
            samp.method_of_sample_procurement AS `Sample Procurement Method`
'@

$filesQuery = @'
MATCH (f:file)-->(parent)
MATCH (f)-[:file_of_sample]->(samp)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (d)-[:diagnosis_of_study_subject]->(ss)
MATCH (tp)-[:tp_of_diagnosis]->(d)
 WHERE ss.disease_subtype IN ["Infiltrating Lobular Carcinoma"]
WITH
        f, parent,p, ss, d,tp, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent,p, ss, d,tp, s, samp,
        f.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
WITH
        f, parent,p, ss, d,tp, s, samp, unit,
        round(factor * value)/factor AS size
RETURN Distinct
    f.file_name AS `File Name`,
    head(labels(samp)) AS `Association`,
    f.file_description AS `Description`,
    f.file_format AS `File Format`,
     CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    p.program_acronym AS `Program Code`,
    s.study_acronym AS `Arm`,
    ss.study_subject_id AS `Case ID`,
    samp.sample_id AS `Sample ID`
    order by f.file_name
'@

# Tab names first (A3, A4), then queries (B3, B4), to match shared-string
# insertion order: SamplesTab, FilesTab, SamplesQuery, FilesQuery
$ws.Range("A3").Value = "SamplesTab"
$ws.Range("A4").Value = "FilesTab"
$ws.Range("B3").Value = $samplesQuery
$ws.Range("B4").Value = $filesQuery

# Row 3: SamplesTab
$ws.Range("C3").Value = $ws.Range("C2").Value2
$ws.Range("D3").Value = $ws.Range("D2").Value2
$ws.Range("E3").Value = $ws.Range("E2").Value2
$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true

# Row 4: FilesTab
$ws.Range("C4").Value = $ws.Range("C2").Value2
$ws.Range("D4").Value = $ws.Range("D2").Value2
$ws.Range("E4").Value = $ws.Range("E2").Value2
$ws.Range("B4").WrapText = $true
$ws.Range("C4").WrapText = $true

$ws.Range("B4").Select()
